$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TCar")

# Replace the #REF! formulas in column D (rows 5-14) with a plain default value,
# and fill the previously-empty H:N columns for the same rows with a default 0.
for ($row = 5; $row -le 14; $row++) {
    $ws.Cells.Item($row, 4).Value = 0   # D

    $ws.Cells.Item($row, 8).Value = 0   # H
    $ws.Cells.Item($row, 9).Value = 0   # I
    $ws.Cells.Item($row, 10).Value = 0  # J
    $ws.Cells.Item($row, 11).Value = 0  # K
    $ws.Cells.Item($row, 12).Value = 0  # L
    $ws.Cells.Item($row, 13).Value = 0  # M
    $ws.Cells.Item($row, 14).Value = 0  # N
}

# Move the selection/viewport: activate TCar, scroll so column A is visible
# again, and select A5 (matches the recorded cursor position after editing).
$ws.Activate()
$ws.Range("A5").Select()
